$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (new A value, new B text)
$updates = @{
    4  = @(25, ' bwd iat max')
    5  = @(26, ' bwd iat min')
    6  = @(13, ' bwd packet length min')
    7  = @(29, ' bwd packets/s')
    9  = @(4,  ' flow duration')
    10 = @(17, ' flow iat max')
    11 = @(15, ' flow iat mean')
    12 = @(18, ' flow iat min')
    13 = @(16, ' flow iat std')
    14 = @(14, ' flow packets/s')
    15 = @(27, ' fwd header length')
    16 = @(32, ' fwd header length.1')
    17 = @(22, ' fwd iat max')
    18 = @(20, ' fwd iat mean')
    19 = @(23, ' fwd iat min')
    20 = @(21, ' fwd iat std')
    21 = @(9,  ' fwd packet length max')
    22 = @(10, ' fwd packet length min')
    23 = @(11, ' fwd packet length std')
    24 = @(39, ' idle max')
    25 = @(38, ' idle std')
    26 = @(3,  ' protocol')
    27 = @(31, ' rst flag count')
    28 = @(1,  ' source port')
    29 = @(30, ' syn flag count')
    30 = @(6,  ' total backward packets')
    31 = @(5,  ' total fwd packets')
    32 = @(8,  ' total length of bwd packets')
    33 = @(34, 'active mean')
    34 = @(24, 'bwd iat total')
    35 = @(12, 'bwd packet length max')
    36 = @(19, 'fwd iat total')
    37 = @(28, 'fwd packets/s')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
}
